$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:G1): drop the bold / thin-border / centered style that
#     was previously applied so the header cells fall back to the default
#     (unstyled) format. ClearFormats() removes the font/border/alignment
#     formatting entirely instead of just pointing at a different style. ---
$ws.Range("A1:G1").ClearFormats()

# --- New "eliminado" (soft-delete / papelera) column ---
$ws.Range("H1").Value = "eliminado"

# Row 2: keep column H present-but-blank for this existing product.
$ws.Range("H2").Borders.LineStyle = -4142  # xlLineStyleNone (no-op touch so the cell is materialized but stays empty)

# Row 3: price/stock correction + new "eliminado" flag (not deleted => 0)
$ws.Range("D3").Value = 61763
$ws.Range("E3").Value = 4
$ws.Range("H3").Value = 0

# Row 4: brand-new product row
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Yuval"
$ws.Range("C4").Value = "kllklk"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "img/producto_3.jpg"
$ws.Range("H4").Value = $false
